# Update the "Notes" column for residence_icb_code (row 2) and
# residence_icb_name (row 3) to clarify how missing LSOA is handled.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newNote = "Raw extract constrained to any midlands resident OR admission to midlands based provider. Last known LSOA used when missing from record."

$ws.Range("D2").Value = $newNote
$ws.Range("D3").Value = $newNote

$ws.Range("D4").Select()
